$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-03 Wednesday" "2024-07-04 Thursday"

Replace-Text "51÷2=" "96÷5="
Replace-Text "93÷3=" "91÷4="
Replace-Text "35÷4=" "61÷4="
Replace-Text "89÷9=" "74÷5="
Replace-Text "77÷2=" "27÷9="
Replace-Text "90÷2=" "55÷4="
Replace-Text "58÷2=" "71÷6="
Replace-Text "70÷4=" "45÷2="
Replace-Text "25÷5=" "55÷6="
Replace-Text "91÷2=" "17÷3="
Replace-Text "43÷8=" "88÷2="
Replace-Text "23÷2=" "88÷8="
Replace-Text "69÷9=" "69÷2="
Replace-Text "88÷9=" "82÷5="
Replace-Text "84÷8=" "89÷4="
Replace-Text "14÷5=" "45÷4="
Replace-Text "28÷5=" "26÷4="
Replace-Text "14÷6=" "69÷5="
Replace-Text "56÷6=" "11÷4="
Replace-Text "52÷7=" "70÷3="
Replace-Text "10÷9=" "86÷2="
Replace-Text "80÷3=" "88÷3="
Replace-Text "27÷3=" "36÷4="
Replace-Text "22÷4=" "37÷6="
Replace-Text "35÷3=" "29÷9="
